{"js": "// Remove the last row of the second table (an empty spacer row) and\n// collapse the two trailing empty \"centered\" paragraphs after that table\n// into a single plain paragraph (no alignment / paragraph formatting).\n\nconst body = context.document.body;\n\n// --- 1) Delete the last row of the second table ---------------------------\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst lastTable = tables.items[tables.items.length - 1];\nconst rows = lastTable.rows;\nrows.load(\"items\");\nawait context.sync();\n\nconst lastRow = rows.items[rows.items.length - 1];\nlastRow.delete();\nawait context.sync();\n\n// --- 2) Collapse the two trailing empty paragraphs into one ---------------\n// body.paragraphs enumerates every paragraph in the story (including the\n// ones nested in table cells), so the very last two entries are the pair of\n// centered empty paragraphs that sit after the table.\nconst paras = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\nconst n = paras.items.length;\n// Word never allows deleting the final paragraph mark of the body, so\n// delete the second-to-last paragraph and keep the last one.\nparas.items[n - 2].delete();\nawait context.sync();\n\n// Re-fetch so we operate on the now-last paragraph and strip its\n// paragraph formatting (the centered alignment) so it becomes a bare,\n// unformatted paragraph.\nconst paras2 = body.paragraphs;\nparas2.load(\"items\");\nawait context.sync();\n\nconst remaining = paras2.items[paras2.items.length - 1];\nremaining.alignment = Word.Alignment.left;\nawait context.sync();\n", "ps1": "# Remove the last row of the second table (an empty spacer row) and\n# collapse the two trailing empty \"centered\" paragraphs after that table\n# into a single plain paragraph (no alignment / paragraph formatting).\n\n$d = $word.ActiveDocument\n\n# --- 1) Delete the last row of the second (last) table ---------------------\n$table = $d.Tables.Item($d.Tables.Count)\n$rowCount = $table.Rows.Count\n$table.Rows.Item($rowCount).Delete()\n\n# --- 2) Collapse the two trailing empty paragraphs into one -----------------\n# Word never allows deleting the final paragraph mark of the body, so delete\n# the second-to-last paragraph (Last.Previous()) and keep the last one.\n$last = $d.Paragraphs.Last\n$prev = $last.Previous()\n$prev.Range.Delete()\n\n# Strip the remaining paragraph's formatting (the centered alignment) so it\n# becomes a bare, unformatted paragraph.\n$d2 = $word.ActiveDocument\n$remaining = $d2.Paragraphs.Last\n$remaining.Format.Alignment = 0\n"}
